# Update the KPI computations on three worksheets: Productdata, Capacity,
# and ProcessingTime, per the commit "Change the computations of the KPIs".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Productdata sheet: columns C (flag) and E (computed KPI ratio) for rows 2-23
# ---------------------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("Productdata")

$wsProduct.Range("C2").Value = 0
$wsProduct.Range("E2").Value = 0.1782

$wsProduct.Range("E3").Value = 0.05827499999999999

$wsProduct.Range("E4").Value = 0.05834444444444443

$wsProduct.Range("C5").Value = 0
$wsProduct.Range("E5").Value = 0.05803611111111111

$wsProduct.Range("E6").Value = 0.05693888888888888

$wsProduct.Range("E7").Value = 0.05706666666666666

$wsProduct.Range("E8").Value = 0.05682222222222223

$wsProduct.Range("E9").Value = 0.06011388888888888

$wsProduct.Range("E10").Value = 0.05999999999999999

$wsProduct.Range("E11").Value = 0.05969166666666666

$wsProduct.Range("E12").Value = 0.06008333333333334

$wsProduct.Range("C13").Value = 0
$wsProduct.Range("E13").Value = 0.1706833333333333

$wsProduct.Range("C14").Value = 0
$wsProduct.Range("E14").Value = 0.05701388888888888

$wsProduct.Range("C15").Value = 0
$wsProduct.Range("E15").Value = 0.05624722222222223

$wsProduct.Range("C16").Value = 0
$wsProduct.Range("E16").Value = 0.05566111111111111

$wsProduct.Range("C17").Value = 0
$wsProduct.Range("E17").Value = 0.05582499999999999

$wsProduct.Range("C18").Value = 0
$wsProduct.Range("E18").Value = 0.05595

$wsProduct.Range("C19").Value = 0
$wsProduct.Range("E19").Value = 0.05572222222222222

$wsProduct.Range("E20").Value = 0.06574999999999999

$wsProduct.Range("E21").Value = 0.0698

$wsProduct.Range("E22").Value = 0.0864

$wsProduct.Range("E23").Value = 0.26585

# ---------------------------------------------------------------------------
# Capacity sheet: column B for rows 2-23
# ---------------------------------------------------------------------------
$wsCapacity = $wb.Worksheets.Item("Capacity")

$wsCapacity.Range("B2").Value = 30
$wsCapacity.Range("B3").Value = 15
$wsCapacity.Range("B4").Value = 15
$wsCapacity.Range("B5").Value = 25
$wsCapacity.Range("B6").Value = 15
$wsCapacity.Range("B7").Value = 5
$wsCapacity.Range("B8").Value = 10
$wsCapacity.Range("B9").Value = 5
$wsCapacity.Range("B10").Value = 15
$wsCapacity.Range("B11").Value = 20
$wsCapacity.Range("B12").Value = 20
$wsCapacity.Range("B13").Value = 30
$wsCapacity.Range("B14").Value = 5
$wsCapacity.Range("B15").Value = 25
$wsCapacity.Range("B16").Value = 15
$wsCapacity.Range("B17").Value = 25
$wsCapacity.Range("B18").Value = 25
$wsCapacity.Range("B19").Value = 15
$wsCapacity.Range("B20").Value = 45
$wsCapacity.Range("B21").Value = 45
$wsCapacity.Range("B22").Value = 45
$wsCapacity.Range("B23").Value = 180

# ---------------------------------------------------------------------------
# ProcessingTime sheet: diagonal cells (one per product/row)
# ---------------------------------------------------------------------------
$wsProcessing = $wb.Worksheets.Item("ProcessingTime")

$wsProcessing.Range("B2").Value = 2
$wsProcessing.Range("C3").Value = 3
$wsProcessing.Range("G7").Value = 1
$wsProcessing.Range("H8").Value = 2
$wsProcessing.Range("I9").Value = 1
$wsProcessing.Range("J10").Value = 3
$wsProcessing.Range("K11").Value = 4
$wsProcessing.Range("L12").Value = 4
$wsProcessing.Range("P16").Value = 3
$wsProcessing.Range("Q17").Value = 5
$wsProcessing.Range("R18").Value = 5
$wsProcessing.Range("S19").Value = 3
$wsProcessing.Range("T20").Value = 1
$wsProcessing.Range("U21").Value = 1
$wsProcessing.Range("V22").Value = 1
$wsProcessing.Range("W23").Value = 4
